$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J (copy header style from H1, then set values)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for columns I and J, rows 2-27
$data = @{
    2  = @(2, 5)
    3  = @(5, 6)
    4  = @(8, 9)
    5  = @(6, 8)
    6  = @(7, 8)
    7  = @(3, 5)
    8  = @(6, 7)
    9  = @(6, 7)
    10 = @(6, 7)
    11 = @(7, 7)
    12 = @(8, 9)
    13 = @(1, 1)
    14 = @(1, 4)
    15 = @(6, 6)
    16 = @(1, 4)
    17 = @(6, 6)
    18 = @(1, 3)
    19 = @(4, 6)
    20 = @(6, 6)
    21 = @(5, 7)
    22 = @(5, 7)
    23 = @(5, 6)
    24 = @(4, 6)
    25 = @(5, 5)
    26 = @(4, 6)
    27 = @(1, 2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
